$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.672619938850403
$ws.Range("B1").Value = 4.89542818069458
$ws.Range("C1").Value = 3.560811758041382
$ws.Range("D1").Value = 1.693607211112976
$ws.Range("E1").Value = 1.026601195335388
